$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.688262462615967
$ws.Range("B1").Value = 2.450698614120483
$ws.Range("C1").Value = 2.1451096534729
$ws.Range("D1").Value = 1.755123615264893
$ws.Range("E1").Value = 1.380583047866821
